# Updates market-price-derived metrics (currentAveragePrice, NQ/HQ prices, Leve profit columns)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, reflecting a scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 39
$ws.Range("H39").Value = 152.2
$ws.Range("I39").Value = 89.25
$ws.Range("K39").Value = 267.75
$ws.Range("M39").Value = 28.25

# Row 42
$ws.Range("H42").Value = 200130.84
$ws.Range("I42").Value = 236494.64
$ws.Range("J42").Value = 130
$ws.Range("K42").Value = 709483.92
$ws.Range("L42").Value = 390
$ws.Range("M42").Value = -709253.92
$ws.Range("N42").Value = -850

# Row 58
$ws.Range("H58").Value = 2081.7144
$ws.Range("J58").Value = 4611
$ws.Range("L58").Value = 13833
$ws.Range("N58").Value = -14133

# Row 86
$ws.Range("H86").Value = 8778262
$ws.Range("J86").Value = 17549524
$ws.Range("L86").Value = 17549524
$ws.Range("N86").Value = -17551770

# Row 89
$ws.Range("H89").Value = 8778262
$ws.Range("J89").Value = 17549524
$ws.Range("L89").Value = 87747620
$ws.Range("N89").Value = -87758852

# Row 92
$ws.Range("H92").Value = 4630851
$ws.Range("I92").Value = 970.25
$ws.Range("K92").Value = 970.25
$ws.Range("M92").Value = 277.75

# Row 100
$ws.Range("H100").Value = 5714.8184
$ws.Range("J100").Value = 6943.643
$ws.Range("L100").Value = 6943.643
$ws.Range("N100").Value = -8025.643

# Row 106
$ws.Range("H106").Value = 2600.3333
$ws.Range("I106").Value = 2825.5
$ws.Range("K106").Value = 2825.5
$ws.Range("M106").Value = -2194.5

# Row 116
$ws.Range("H116").Value = 5584.1665
$ws.Range("I116").Value = 4626.25
$ws.Range("K116").Value = 4626.25
$ws.Range("M116").Value = -1184.25

# Row 125
$ws.Range("H125").Value = 2001
$ws.Range("I125").Value = 2141
$ws.Range("J125").Value = 1961
$ws.Range("K125").Value = 19269
$ws.Range("L125").Value = 17649
$ws.Range("M125").Value = -16809
$ws.Range("N125").Value = -22569

# Row 141
$ws.Range("H141").Value = 2779.5881
$ws.Range("I141").Value = 2779.5881
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 8338.764299999999
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -3158.764299999999
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 11882.145
$ws.Range("I32").Value = 8854.825000000001
$ws.Range("K32").Value = 8854.825000000001
$ws.Range("M32").Value = -8567.825000000001

# Row 74
$ws.Range("H74").Value = 2115.25
$ws.Range("J74").Value = 3823.5
$ws.Range("L74").Value = 3823.5
$ws.Range("N74").Value = -5571.5

# Row 77
$ws.Range("H77").Value = 2115.25
$ws.Range("J77").Value = 3823.5
$ws.Range("L77").Value = 19117.5
$ws.Range("N77").Value = -27853.5

# Row 97
$ws.Range("H97").Value = 1033.2142
$ws.Range("I97").Value = 1073.6818
$ws.Range("J97").Value = 884.8333
$ws.Range("K97").Value = 1073.6818
$ws.Range("L97").Value = 884.8333
$ws.Range("M97").Value = -577.6818000000001
$ws.Range("N97").Value = -1876.8333

# Row 122
$ws.Range("H122").Value = 3809.7307
$ws.Range("I122").Value = 2797.5789
$ws.Range("K122").Value = 8392.736699999999
$ws.Range("M122").Value = -5942.736699999999

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 13045334
$ws.Range("I94").Value = 5001183.5
$ws.Range("K94").Value = 5001183.5
$ws.Range("M94").Value = -5000732.5

# Row 130
$ws.Range("H130").Value = 51796.668
$ws.Range("J130").Value = 51796.668
$ws.Range("L130").Value = 51796.668
$ws.Range("N130").Value = -61836.668

$ws = $wb.Worksheets.Item("CRP")
# Row 3
$ws.Range("H3").Value = 5666
$ws.Range("I3").Value = 5666
$ws.Range("K3").Value = 5666
$ws.Range("M3").Value = -5553

# Row 22
$ws.Range("H22").Value = 1160.1
$ws.Range("I22").Value = 821.2
$ws.Range("K22").Value = 821.2
$ws.Range("M22").Value = -471.2

# Row 53
$ws.Range("H53").Value = 108332.664
$ws.Range("J53").Value = 108332.664
$ws.Range("L53").Value = 108332.664
$ws.Range("N53").Value = -109546.664

# Row 60
$ws.Range("H60").Value = 8999.6
$ws.Range("I60").Value = 8999.6
$ws.Range("K60").Value = 8999.6
$ws.Range("M60").Value = -8488.6

# Row 134
$ws.Range("H134").Value = 3558.92
$ws.Range("I134").Value = 2540.6843
$ws.Range("K134").Value = 7622.0529
$ws.Range("M134").Value = -5087.0529

$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 204.53847
$ws.Range("I14").Value = 204.53847
$ws.Range("K14").Value = 613.61541
$ws.Range("M14").Value = -440.61541

# Row 37
$ws.Range("H37").Value = 235675.64
$ws.Range("J37").Value = 235675.64
$ws.Range("L37").Value = 707026.92
$ws.Range("N37").Value = -707250.92

# Row 133
$ws.Range("H133").Value = 2944.5
$ws.Range("I133").Value = 2944.5
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 8833.5
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -3773.5
$ws.Range("N133").ClearContents()

# Row 134
$ws.Range("H134").Value = 6710.7144
$ws.Range("I134").Value = 2378.8462
$ws.Range("K134").Value = 7136.5386
$ws.Range("M134").Value = -2066.5386

# Row 139
$ws.Range("H139").Value = 7010.657
$ws.Range("I139").Value = 3947.75
$ws.Range("J139").Value = 8608.695
$ws.Range("K139").Value = 11843.25
$ws.Range("L139").Value = 25826.085
$ws.Range("M139").Value = -6703.25
$ws.Range("N139").Value = -36106.085

# Row 140
$ws.Range("H140").Value = 1967.8334
$ws.Range("I140").Value = 952.75
$ws.Range("K140").Value = 2858.25
$ws.Range("M140").Value = 2321.75

# Row 141
$ws.Range("H141").Value = 19190.477
$ws.Range("I141").Value = 3000
$ws.Range("K141").Value = 9000
$ws.Range("M141").Value = -3820

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 41814.613
$ws.Range("I70").Value = 61101
$ws.Range("J70").Value = 6748.4546
$ws.Range("K70").Value = 61101
$ws.Range("L70").Value = 6748.4546
$ws.Range("M70").Value = -60831
$ws.Range("N70").Value = -7288.4546

# Row 73
$ws.Range("H73").Value = 41814.613
$ws.Range("I73").Value = 61101
$ws.Range("J73").Value = 6748.4546
$ws.Range("K73").Value = 61101
$ws.Range("L73").Value = 6748.4546
$ws.Range("M73").Value = -60165
$ws.Range("N73").Value = -8620.454600000001

# Row 80
$ws.Range("H80").Value = 19676794
$ws.Range("I80").Value = 95698.414
$ws.Range("J80").Value = 66671424
$ws.Range("K80").Value = 95698.414
$ws.Range("L80").Value = 66671424
$ws.Range("M80").Value = -94700.414
$ws.Range("N80").Value = -66673420

# Row 83
$ws.Range("H83").Value = 19676794
$ws.Range("I83").Value = 95698.414
$ws.Range("J83").Value = 66671424
$ws.Range("K83").Value = 478492.07
$ws.Range("L83").Value = 333357120
$ws.Range("M83").Value = -473500.07
$ws.Range("N83").Value = -333367104

# Row 102
$ws.Range("H102").Value = 2304.7837
$ws.Range("I102").Value = 1751
$ws.Range("K102").Value = 1751
$ws.Range("M102").Value = -129

# Row 107
$ws.Range("H107").Value = 537
$ws.Range("I107").Value = 537
$ws.Range("K107").Value = 537
$ws.Range("M107").Value = 1383

# Row 122
$ws.Range("H122").Value = 4974.6113
$ws.Range("I122").Value = 2747.111
$ws.Range("J122").Value = 7202.1113
$ws.Range("K122").Value = 8241.332999999999
$ws.Range("L122").Value = 21606.3339
$ws.Range("M122").Value = -5791.332999999999
$ws.Range("N122").Value = -26506.3339

# Row 132
$ws.Range("H132").Value = 3263.8718
$ws.Range("I132").Value = 3263.8718
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9791.615399999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7261.615399999999
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4343.7617
$ws.Range("I7").Value = 2696.818
$ws.Range("J7").Value = 6155.4
$ws.Range("K7").Value = 2696.818
$ws.Range("L7").Value = 6155.4
$ws.Range("M7").Value = -2584.818
$ws.Range("N7").Value = -6379.4

# Row 46
$ws.Range("H46").Value = 8318.892
$ws.Range("I46").Value = 1166.5
$ws.Range("K46").Value = 1166.5
$ws.Range("M46").Value = -978.5

# Row 61
$ws.Range("H61").Value = 2603.2646
$ws.Range("I61").Value = 1682.68
$ws.Range("J61").Value = 5160.4443
$ws.Range("K61").Value = 1682.68
$ws.Range("L61").Value = 5160.4443
$ws.Range("M61").Value = -1480.68
$ws.Range("N61").Value = -5564.4443

# Row 113
$ws.Range("H113").Value = 2603.2646
$ws.Range("I113").Value = 1682.68
$ws.Range("J113").Value = 5160.4443
$ws.Range("K113").Value = 1682.68
$ws.Range("L113").Value = 5160.4443
$ws.Range("M113").Value = 487.3199999999999
$ws.Range("N113").Value = -9500.444299999999

# Row 126
$ws.Range("H126").Value = 4343.7617
$ws.Range("I126").Value = 2696.818
$ws.Range("J126").Value = 6155.4
$ws.Range("K126").Value = 8090.454000000001
$ws.Range("L126").Value = 18466.2
$ws.Range("M126").Value = -5620.454000000001
$ws.Range("N126").Value = -23406.2

$ws = $wb.Worksheets.Item("WVR")
# Row 8
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()

# Row 28
$ws.Range("H28").Value = 17250
$ws.Range("I28").Value = 14500
$ws.Range("J28").Value = 20000
$ws.Range("K28").Value = 14500
$ws.Range("L28").Value = 20000
$ws.Range("M28").Value = -14152
$ws.Range("N28").Value = -20696

# Row 49
$ws.Range("H49").Value = 32270.273
$ws.Range("I49").Value = 24998
$ws.Range("K49").Value = 24998
$ws.Range("M49").Value = -24768

# Row 81
$ws.Range("H81").Value = 66671770
$ws.Range("I81").Value = 3001
$ws.Range("J81").Value = 83338960
$ws.Range("K81").Value = 6002
$ws.Range("L81").Value = 166677920
$ws.Range("M81").Value = -4941
$ws.Range("N81").Value = -166680042

# Row 84
$ws.Range("H84").Value = 66671770
$ws.Range("I84").Value = 3001
$ws.Range("J84").Value = 83338960
$ws.Range("K84").Value = 30010
$ws.Range("L84").Value = 833389600
$ws.Range("M84").Value = -24706
$ws.Range("N84").Value = -833400208
